$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 4579
$ws.Range("I4").Value = 5844.143
$ws.Range("K4").Value = 5844.143
$ws.Range("M4").Value = -5730.143

$ws.Range("H18").Value = 7642.7144
$ws.Range("I18").Value = 8791.5
$ws.Range("J18").Value = 750
$ws.Range("K18").Value = 8791.5
$ws.Range("L18").Value = 750
$ws.Range("M18").Value = -8507.5
$ws.Range("N18").Value = -1318

$ws.Range("H33").Value = 19262.354
$ws.Range("I33").Value = 29592.637
$ws.Range("J33").Value = 323.5
$ws.Range("K33").Value = 29592.637
$ws.Range("L33").Value = 323.5
$ws.Range("M33").Value = -29363.637
$ws.Range("N33").Value = -781.5

$ws.Range("H38").Value = 376.63635
$ws.Range("I38").Value = 114.3
$ws.Range("K38").Value = 342.9
$ws.Range("M38").Value = 29.10000000000002

$ws.Range("J76").Value = 7500
$ws.Range("L76").Value = 7500
$ws.Range("N76").Value = -8130

$ws.Range("J79").Value = 7500
$ws.Range("L79").Value = 7500
$ws.Range("N79").Value = -9684

$ws.Range("H107").Value = 725.0526
$ws.Range("I107").Value = 725.0526
$ws.Range("K107").Value = 725.0526
$ws.Range("M107").Value = 1194.9474

$ws.Range("H111").Value = 2029.6666
$ws.Range("I111").Value = 1794.5
$ws.Range("K111").Value = 5383.5
$ws.Range("M111").Value = -2316.5

$ws.Range("H112").Value = 2220.8
$ws.Range("J112").Value = 2285.6155
$ws.Range("L112").Value = 6856.8465
$ws.Range("N112").Value = -9072.8465

$ws.Range("H135").Value = 1902.75
$ws.Range("I135").Value = 2046
$ws.Range("K135").Value = 18414
$ws.Range("M135").Value = -15879

$ws.Range("H137").Value = 1738.6428
$ws.Range("I137").Value = 1624.7142
$ws.Range("K137").Value = 4874.142599999999
$ws.Range("M137").Value = -2324.142599999999

$ws.Range("H141").Value = 8099.7896
$ws.Range("I141").Value = 6806.3125
$ws.Range("K141").Value = 20418.9375
$ws.Range("M141").Value = -15238.9375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3736.5
$ws.Range("I45").Value = 2661.303
$ws.Range("J45").Value = 5823.647
$ws.Range("K45").Value = 2661.303
$ws.Range("L45").Value = 5823.647
$ws.Range("M45").Value = -2284.303
$ws.Range("N45").Value = -6577.647

$ws.Range("H102").Value = 4854.9
$ws.Range("I102").Value = 3873.6
$ws.Range("K102").Value = 3873.6
$ws.Range("M102").Value = -2251.6

$ws.Range("H122").Value = 1811.1052
$ws.Range("I122").Value = 1871.2941
$ws.Range("K122").Value = 5613.8823
$ws.Range("M122").Value = -3163.8823

$ws.Range("H132").Value = 3915.2063
$ws.Range("I132").Value = 3228.6875
$ws.Range("J132").Value = 6112.067
$ws.Range("K132").Value = 9686.0625
$ws.Range("L132").Value = 18336.201
$ws.Range("M132").Value = -7156.0625
$ws.Range("N132").Value = -23396.201

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 6716.3335
$ws.Range("I22").Value = 6716.3335
$ws.Range("K22").Value = 6716.3335
$ws.Range("M22").Value = -6543.3335

$ws.Range("H40").Value = 10448
$ws.Range("J40").Value = 10448
$ws.Range("L40").Value = 10448
$ws.Range("N40").Value = -10978

$ws.Range("H96").Value = 12125.2
$ws.Range("I96").Value = 13788.75
$ws.Range("J96").Value = 5471
$ws.Range("K96").Value = 13788.75
$ws.Range("L96").Value = 5471
$ws.Range("M96").Value = -11042.75
$ws.Range("N96").Value = -10963

$ws.Range("H99").Value = 4732.448
$ws.Range("I99").Value = 3730.524
$ws.Range("J99").Value = 7362.5
$ws.Range("K99").Value = 3730.524
$ws.Range("L99").Value = 7362.5
$ws.Range("M99").Value = -2232.524
$ws.Range("N99").Value = -10358.5

$ws.Range("H105").Value = 4212.1665
$ws.Range("I105").Value = 4187.6875
$ws.Range("K105").Value = 4187.6875
$ws.Range("M105").Value = -2440.6875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9061
$ws.Range("I31").Value = 7459
$ws.Range("K31").Value = 7459
$ws.Range("M31").Value = -7164

$ws.Range("H34").Value = 9061
$ws.Range("I34").Value = 7459
$ws.Range("K34").Value = 7459
$ws.Range("M34").Value = -7257

$ws.Range("H107").Value = 1904.3334
$ws.Range("I107").Value = 1207.5834
$ws.Range("K107").Value = 1207.5834
$ws.Range("M107").Value = 712.4166

$ws.Range("H132").Value = 3075
$ws.Range("J132").Value = 3075
$ws.Range("L132").Value = 9225
$ws.Range("N132").Value = -14285

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 10220128
$ws.Range("I4").Value = 5231184.5
$ws.Range("K4").Value = 15693553.5
$ws.Range("M4").Value = -15693441.5

$ws.Range("H14").Value = 416.6
$ws.Range("I14").Value = 416.6
$ws.Range("K14").Value = 1249.8
$ws.Range("M14").Value = -1076.8

$ws.Range("H33").Value = 743.0833
$ws.Range("I33").Value = 124
$ws.Range("J33").Value = 1609.8
$ws.Range("K33").Value = 744
$ws.Range("L33").Value = 9658.799999999999
$ws.Range("M33").Value = -461
$ws.Range("N33").Value = -10224.8

$ws.Range("H37").Value = 56557.48
$ws.Range("J37").Value = 56557.48
$ws.Range("L37").Value = 169672.44
$ws.Range("N37").Value = -169896.44

$ws.Range("H68").Value = 2618.7778
$ws.Range("J68").Value = 3028.1667
$ws.Range("L68").Value = 9084.500100000001
$ws.Range("N68").Value = -10706.5001

$ws.Range("H71").Value = 2618.7778
$ws.Range("J71").Value = 3028.1667
$ws.Range("L71").Value = 27253.5003
$ws.Range("N71").Value = -35365.5003

$ws.Range("H86").Value = 564.1429000000001
$ws.Range("I86").Value = 450
$ws.Range("J86").Value = 609.8
$ws.Range("K86").Value = 1350
$ws.Range("L86").Value = 1829.4
$ws.Range("M86").Value = -164
$ws.Range("N86").Value = -4201.4

$ws.Range("H89").Value = 564.1429000000001
$ws.Range("I89").Value = 450
$ws.Range("J89").Value = 609.8
$ws.Range("K89").Value = 4050
$ws.Range("L89").Value = 5488.2
$ws.Range("M89").Value = 1878
$ws.Range("N89").Value = -17344.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 12930.375
$ws.Range("I70").Value = 4499
$ws.Range("J70").Value = 14134.857
$ws.Range("K70").Value = 4499
$ws.Range("L70").Value = 14134.857
$ws.Range("M70").Value = -4229
$ws.Range("N70").Value = -14674.857

$ws.Range("H73").Value = 12930.375
$ws.Range("I73").Value = 4499
$ws.Range("J73").Value = 14134.857
$ws.Range("K73").Value = 4499
$ws.Range("L73").Value = 14134.857
$ws.Range("M73").Value = -3563
$ws.Range("N73").Value = -16006.857

$ws.Range("H97").Value = 2266.4546
$ws.Range("I97").Value = 690.9167
$ws.Range("K97").Value = 690.9167
$ws.Range("M97").Value = -194.9167

$ws.Range("H111").Value = 50000
$ws.Range("J111").Value = 50000
$ws.Range("L111").Value = 50000
$ws.Range("N111").Value = -56134

$ws.Range("H122").Value = 2884.2222
$ws.Range("I122").Value = 2656.9375
$ws.Range("K122").Value = 7970.8125
$ws.Range("M122").Value = -5520.8125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4842.6665
$ws.Range("I40").Value = 3917.6
$ws.Range("K40").Value = 3917.6
$ws.Range("M40").Value = -3781.6

$ws.Range("H100").Value = 5004571.5
$ws.Range("I100").Value = 12503929
$ws.Range("J100").Value = 5000
$ws.Range("K100").Value = 12503929
$ws.Range("L100").Value = 5000
$ws.Range("M100").Value = -12503388
$ws.Range("N100").Value = -6082

$ws.Range("H122").Value = 4933.0835
$ws.Range("I122").Value = 3642.4285
$ws.Range("K122").Value = 10927.2855
$ws.Range("M122").Value = -8477.2855

$ws.Range("H132").Value = 13037.929
$ws.Range("I132").Value = 13664.846
$ws.Range("K132").Value = 40994.538
$ws.Range("M132").Value = -38464.538

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()

$ws.Range("H81").Value = 3203.2856
$ws.Range("I81").Value = 1549.8
$ws.Range("K81").Value = 3099.6
$ws.Range("M81").Value = -2038.6

$ws.Range("H84").Value = 3203.2856
$ws.Range("I84").Value = 1549.8
$ws.Range("K84").Value = 15498
$ws.Range("M84").Value = -10194

$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()

$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

$ws.Range("H113").Value = 341.09525
$ws.Range("I113").Value = 362.6316
$ws.Range("J113").Value = 136.5
$ws.Range("K113").Value = 1087.8948
$ws.Range("L113").Value = 409.5
$ws.Range("M113").Value = 1082.1052
$ws.Range("N113").Value = -4749.5

$ws.Range("H122").Value = 2679.1292
$ws.Range("I122").Value = 2074.2173
$ws.Range("J122").Value = 4418.25
$ws.Range("K122").Value = 6222.651899999999
$ws.Range("L122").Value = 13254.75
$ws.Range("M122").Value = -3772.651899999999
$ws.Range("N122").Value = -18154.75

$ws.Range("H133").Value = 200000
$ws.Range("J133").Value = 200000
$ws.Range("L133").Value = 200000
$ws.Range("N133").Value = -210120
